$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 11 (existing data, and formatting such as the
# date-style on column D, shifts down automatically with the rest of the rows).
$ws.Rows("11:14").Insert()

# New week of "Frutilla" price data (2023-08-07 / serial 45145), one row per
# Calidad grade, mirroring the layout of the rest of the sheet.
$rows = @(
    @{ Row = 11; Calidad = "Especial"; Volumen = 100; Min = 7000; Max = 8000; Prom = 7600; KgPrecio = 2533 },
    @{ Row = 12; Calidad = "Primera";  Volumen = 140; Min = 5000; Max = 6000; Prom = 5500; KgPrecio = 1833 },
    @{ Row = 13; Calidad = "Segunda";  Volumen = 190; Min = 4000; Max = 5000; Prom = 4474; KgPrecio = 1491 },
    @{ Row = 14; Calidad = "Tercera";  Volumen = 160; Min = 3000; Max = 4000; Prom = 3625; KgPrecio = 1208 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 1
    $ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($row, 3).Value = "Arica y Parinacota"
    $ws.Cells.Item($row, 4).Value = 45145
    $ws.Cells.Item($row, 5).Value = 15
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100101
    $ws.Cells.Item($row, 8).Value = "Berries"
    $ws.Cells.Item($row, 9).Value = 100112025
    $ws.Cells.Item($row, 10).Value = "Frutilla"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = "`$/bandeja 3 kilos"
    $ws.Cells.Item($row, 18).Value = "Región de Arica y Parinacota"
    $ws.Cells.Item($row, 19).Value = $r.KgPrecio
    $ws.Cells.Item($row, 20).Value = 3
}
